# Updates the cryptos list "Price" (D) and "Volume(1h)" (E) columns with the
# latest refreshed values. These columns are stored as plain text, so for the
# "Price" column (which often looks like a number, e.g. "42.21") we briefly
# force a text NumberFormat before assigning the value; otherwise Excel would
# auto-coerce the text into a numeric value. ClearFormats() afterwards removes
# the temporary formatting so the cell keeps its original (unstyled) look.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.062.42"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.95%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.780.45"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.47"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.03%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4503"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3571"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.42%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07463"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.21"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("E11").Value = "  +2.92%  "
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.05"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.073"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("E15").Value = "  +2.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.778.99"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.79"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06440"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9998"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.19"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.815"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.54%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.074.64"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.37"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.123"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.32%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.91"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.75%  "
$ws.Range("E27").Value = "  +2.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.982.91"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.168"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +6.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.26"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.115"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +6.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.721"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09226"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.687"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.93%  "
$ws.Range("E35").Value = "  +2.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06228"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.23%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02298"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2117"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.020"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6354"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.186"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.397"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.64%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.938"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.89%  "
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.757"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5937"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +2.41%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.08"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.968"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.146"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06902"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.19"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.39%  "
